$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 8207
$ws1.Range("F9").Value  = 7153
$ws1.Range("F14").Value = 712
$ws1.Range("F19").Value = 75
$ws1.Range("F21").Value = 96
$ws1.Range("F22").Value = 11739
$ws1.Range("F25").Value = 2315
$ws1.Range("F27").Value = 3247
$ws1.Range("F29").Value = 2735
$ws1.Range("F32").Value = 1572
$ws1.Range("F38").Value = 5863
$ws1.Range("F43").Value = 163
$ws1.Range("F49").Value = 1130

# Sheet "演出" (Worksheets index 2 / sheet2.xml)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 10
$ws2.Range("F24").Value = 6

# Sheet "全部类型" (Worksheets index 4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 8207
$ws4.Range("F11").Value = 7153
$ws4.Range("F12").Value = 7153
$ws4.Range("F16").Value = 712
$ws4.Range("F22").Value = 96
$ws4.Range("F25").Value = 11739
$ws4.Range("F29").Value = 2315
$ws4.Range("F30").Value = 2315
$ws4.Range("F31").Value = 3247
$ws4.Range("F32").Value = 2735
$ws4.Range("F34").Value = 1573
$ws4.Range("F41").Value = 5863
$ws4.Range("F51").Value = 1130
